$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header labels: switch from slug-like codes to human-readable
# capitalized Spanish labels, and reorder which column holds which label.
$ws.Range("A1").Value = "Horas trabajadas"
$ws.Range("B1").Value = "Personas residentes viviendas familiares"
$ws.Range("C1").Value = "Provincia código"
$ws.Range("D1").Value = "Aragón"
$ws.Range("E1").Value = "Municipio código"
$ws.Range("F1").Value = "Provincia nombre"
$ws.Range("G1").Value = "Sexo"
$ws.Range("H1").Value = "Municipio nombre"

# Row 2 - DSD concept reference for each column
$ws.Range("A2").Value = "iaest-measure:horas-trabajadas"
$ws.Range("B2").Value = "iaest-measure:personas-residentes-viviendas-familiares"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "null"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("G2").Value = "iaest-measure:sexo"
$ws.Range("H2").Value = "sdmx-dimension:refArea"

# Row 3 - component type (measure / dimension)
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "dim"
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "dim"
$ws.Range("G3").Value = "medida"
$ws.Range("H3").Value = "dim"

# Row 4 - datatype / representation URI
$ws.Range("A4").Value = "xsd:string"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "URI-Comunidad"
$ws.Range("E4").Value = "null"
$ws.Range("F4").Value = "URI-Provincia"
$ws.Range("G4").Value = "xsd:string"
$ws.Range("H4").Value = "URI-Municipio"
